$wb = $excel.ActiveWorkbook

# ALC row 28
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 777.4375
$ws.Range("I28").Value = 822.6
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 822.6
$ws.Range("L28").Value = 100
$ws.Range("M28").Value = -337.6
$ws.Range("N28").Value = -1070

# ALC row 33
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 387
$ws.Range("I33").Value = 332.33334
$ws.Range("J33").Value = 551
$ws.Range("K33").Value = 332.33334
$ws.Range("L33").Value = 551
$ws.Range("M33").Value = -103.33334
$ws.Range("N33").Value = -1009

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 332892.72
$ws.Range("I132").Value = 123502.06
$ws.Range("J132").Value = 1002942.8
$ws.Range("K132").Value = 370506.18
$ws.Range("L132").Value = 3008828.4
$ws.Range("M132").Value = -367976.18
$ws.Range("N132").Value = -3013888.4

# ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 54449.688
$ws.Range("J133").Value = 54449.688
$ws.Range("L133").Value = 54449.688
$ws.Range("N133").Value = -64569.688

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3151.8167
$ws.Range("J138").Value = 3890.825
$ws.Range("L138").Value = 11672.475
$ws.Range("N138").Value = -21952.475

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3324.6924
$ws.Range("I122").Value = 2979.484
$ws.Range("K122").Value = 8938.451999999999
$ws.Range("M122").Value = -6488.451999999999

# ARM row 137
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H137").Value = 44531.5
$ws.Range("J137").Value = 44531.5
$ws.Range("L137").Value = 44531.5
$ws.Range("N137").Value = -54731.5

# BSM row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3194.6345
$ws.Range("I134").Value = 1170.2903
$ws.Range("J134").Value = 6182.952
$ws.Range("K134").Value = 3510.8709
$ws.Range("L134").Value = 18548.856
$ws.Range("M134").Value = -975.8708999999999
$ws.Range("N134").Value = -23618.856

# CRP row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2126.0317
$ws.Range("I31").Value = 919.14813
$ws.Range("J31").Value = 3031.1943
$ws.Range("K31").Value = 919.14813
$ws.Range("L31").Value = 3031.1943
$ws.Range("M31").Value = -624.14813
$ws.Range("N31").Value = -3621.1943

# CRP row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2126.0317
$ws.Range("I34").Value = 919.14813
$ws.Range("J34").Value = 3031.1943
$ws.Range("K34").Value = 919.14813
$ws.Range("L34").Value = 3031.1943
$ws.Range("M34").Value = -717.14813
$ws.Range("N34").Value = -3435.1943

# CRP row 98
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1467.2593
$ws.Range("I134").Value = 944.64
$ws.Range("K134").Value = 2833.92
$ws.Range("M134").Value = -298.9200000000001

# CRP row 139
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H139").Value = 38916
$ws.Range("J139").Value = 38916
$ws.Range("L139").Value = 38916
$ws.Range("N139").Value = -49196

# CUL row 34
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 20376.182
$ws.Range("I34").Value = 23674
$ws.Range("J34").Value = 17628
$ws.Range("K34").Value = 71022
$ws.Range("L34").Value = 52884
$ws.Range("M34").Value = -70938
$ws.Range("N34").Value = -53052

# CUL row 58
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 4999.7144
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 4999.7144
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 14999.1432
$ws.Range("M58").Value = $null
$ws.Range("N58").Value = -15255.1432

# CUL row 68
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1575.8889
$ws.Range("I68").Value = 1321
$ws.Range("K68").Value = 3963
$ws.Range("M68").Value = -3152

# CUL row 71
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H71").Value = 1575.8889
$ws.Range("I71").Value = 1321
$ws.Range("K71").Value = 11889
$ws.Range("M71").Value = -7833

# CUL row 113
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 5682444
$ws.Range("I113").Value = 671.63635
$ws.Range("K113").Value = 2014.90905
$ws.Range("M113").Value = 155.09095

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 879.03
$ws.Range("I131").Value = 350
$ws.Range("J131").Value = 912.79785
$ws.Range("K131").Value = 1050
$ws.Range("L131").Value = 2738.39355
$ws.Range("M131").Value = 3990
$ws.Range("N131").Value = -12818.39355

# GSM row 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 33509.09
$ws.Range("J15").Value = 33509.09
$ws.Range("L15").Value = 33509.09
$ws.Range("N15").Value = -34085.09

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 13891949
$ws.Range("I80").Value = 20836540
$ws.Range("J80").Value = 2767.6667
$ws.Range("K80").Value = 20836540
$ws.Range("L80").Value = 2767.6667
$ws.Range("M80").Value = -20835542
$ws.Range("N80").Value = -4763.6667

# GSM row 81
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 33509.09
$ws.Range("J81").Value = 33509.09
$ws.Range("L81").Value = 33509.09
$ws.Range("N81").Value = -35505.09

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 13891949
$ws.Range("I83").Value = 20836540
$ws.Range("J83").Value = 2767.6667
$ws.Range("K83").Value = 104182700
$ws.Range("L83").Value = 13838.3335
$ws.Range("M83").Value = -104177708
$ws.Range("N83").Value = -23822.3335

# GSM row 84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 33509.09
$ws.Range("J84").Value = 33509.09
$ws.Range("L84").Value = 100527.27
$ws.Range("N84").Value = -110511.27

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3420.375
$ws.Range("I102").Value = 2346.8948
$ws.Range("K102").Value = 2346.8948
$ws.Range("M102").Value = -724.8948

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3685.2666
$ws.Range("I132").Value = 2144.0667
$ws.Range("K132").Value = 6432.2001
$ws.Range("M132").Value = -3902.2001

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2892.6553
$ws.Range("I7").Value = 1243.7222
$ws.Range("J7").Value = 5590.909
$ws.Range("K7").Value = 1243.7222
$ws.Range("L7").Value = 5590.909
$ws.Range("M7").Value = -1131.7222
$ws.Range("N7").Value = -5814.909

# LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 19309232
$ws.Range("I22").Value = 25100880
$ws.Range("J22").Value = 3733.3333
$ws.Range("K22").Value = 25100880
$ws.Range("L22").Value = 3733.3333
$ws.Range("M22").Value = -25100585
$ws.Range("N22").Value = -4323.3333

# LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 19309232
$ws.Range("I27").Value = 25100880
$ws.Range("J27").Value = 3733.3333
$ws.Range("K27").Value = 25100880
$ws.Range("L27").Value = 3733.3333
$ws.Range("M27").Value = -25100773
$ws.Range("N27").Value = -3947.3333

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 2892.6553
$ws.Range("I126").Value = 1243.7222
$ws.Range("J126").Value = 5590.909
$ws.Range("K126").Value = 3731.1666
$ws.Range("L126").Value = 16772.727
$ws.Range("M126").Value = -1261.1666
$ws.Range("N126").Value = -21712.727

# WVR row 94
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 29950
$ws.Range("J94").Value = 29950
$ws.Range("L94").Value = 29950
$ws.Range("N94").Value = -31752

# WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4463.24
$ws.Range("I122").Value = 2608.6667
$ws.Range("J122").Value = 5506.4375
$ws.Range("K122").Value = 7826.000100000001
$ws.Range("L122").Value = 16519.3125
$ws.Range("M122").Value = -5376.000100000001
$ws.Range("N122").Value = -21419.3125

# WVR row 135
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 102029.2
$ws.Range("J135").Value = 102029.2
$ws.Range("L135").Value = 102029.2
$ws.Range("N135").Value = -112169.2
